# Update "Case Alerts - Public Exposure Sites" workbook:
# replace the existing exposure-site rows with the new Keysborough entry
# (an "old" record and a "new" superseding record), as uploaded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old data rows (rows 4-8), keeping the header (row 1)
# and the first two data rows (rows 2-3) which will be overwritten below.
$ws.Rows("4:8").Delete()

# Row 2: original ("old") exposure window for the Keysborough site
$ws.Range("A2").Value = "Keysborough"
$ws.Range("B2").Value = "Sikh Temple Keysborough, 198-206 Perry Road, Keysborough"
$ws.Range("C2").Value = "1/01/21 3:00pm-5:00pm"
$ws.Range("D2").Value = "Case visited venue"
$ws.Range("E2").Value = "old"

# Row 3: updated ("new") exposure window for the same Keysborough site
$ws.Range("A3").Value = "Keysborough"
$ws.Range("B3").Value = "Sikh Temple Keysborough, 198-206 Perry Road, Keysborough"
$ws.Range("C3").Value = "1/01/21 3:00pm-6:00pm"
$ws.Range("D3").Value = "Case visited venue"
$ws.Range("E3").Value = "new"

# Resize the edited columns to fit their new content
$ws.Columns("A:D").AutoFit()

# Match the saved selection/active cell
$ws.Range("C3").Select()
